$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.744.16"
$ws.Range("E2").Value = "  -5.87%  "
$ws.Range("D3").Value = "'3.581.46"
$ws.Range("E3").Value = "  -3.45%  "
$ws.Range("E4").Value = "  -0.75%  "
$ws.Range("D5").Value = "'399.54"
$ws.Range("E5").Value = "  -6.42%  "
$ws.Range("D6").Value = "'123.72"
$ws.Range("E6").Value = "  -5.49%  "
$ws.Range("D7").Value = "'3.580.04"
$ws.Range("E7").Value = "  -3.23%  "
$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  -8.61%  "
$ws.Range("D10").Value = "'0.686"
$ws.Range("E10").Value = "  -10.64%  "
$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -15.47%  "
$ws.Range("D12").Value = "'0.0000331"
$ws.Range("E12").Value = "  -10.41%  "
$ws.Range("D13").Value = "'39.16"
$ws.Range("E13").Value = "  -8.30%  "
$ws.Range("D14").Value = "'4.132.70"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "'9.26"
$ws.Range("E15").Value = "  -7.51%  "
$ws.Range("B16").Value = "Uniswap"
$ws.Range("C16").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D16").Value = "'14.09"
$ws.Range("E16").Value = "  +9.48%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.136"
$ws.Range("E17").Value = "  -2.72%  "
$ws.Range("D18").Value = "'3.575.71"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "'18.79"
$ws.Range("E19").Value = "  -9.03%  "
$ws.Range("D20").Value = "'63.775.68"
$ws.Range("E20").Value = "  -5.86%  "
$ws.Range("E21").Value = "  -9.50%  "
$ws.Range("D22").Value = "'396.29"
$ws.Range("E22").Value = "  -11.89%  "
$ws.Range("D23").Value = "'13.97"
$ws.Range("E23").Value = "  -8.01%  "
$ws.Range("D24").Value = "'82.32"
$ws.Range("E24").Value = "  -8.07%  "
$ws.Range("D25").Value = "'2.94"
$ws.Range("E25").Value = "  -6.30%  "
$ws.Range("D26").Value = "'5.44"
$ws.Range("E26").Value = "  +9.29%  "
$ws.Range("D27").Value = "'34.22"
$ws.Range("E27").Value = "  -10.51%  "
$ws.Range("D28").Value = "'3.03"
$ws.Range("E28").Value = "  -8.89%  "
$ws.Range("D29").Value = "'8.80"
$ws.Range("E29").Value = "  -15.29%  "
$ws.Range("D30").Value = "'12.01"
$ws.Range("E30").Value = "  -3.96%  "
$ws.Range("D31").Value = "'2.67"
$ws.Range("E31").Value = "  -4.54%  "
$ws.Range("D32").Value = "'0.113"
$ws.Range("E32").Value = "  -6.72%  "
$ws.Range("D33").Value = "'6.88"
$ws.Range("E33").Value = "  -4.60%  "
$ws.Range("D34").Value = "'0.149"
$ws.Range("E34").Value = "  -6.92%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").Value = "'36.83"
$ws.Range("E36").Value = "  -8.94%  "
$ws.Range("D37").Value = "'54.26"
$ws.Range("E37").Value = "  -3.99%  "
$ws.Range("D38").Value = "'0.0439"
$ws.Range("E38").Value = "  -10.61%  "
$ws.Range("D39").Value = "'0.996"
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").Value = "'0.0₃0657"
$ws.Range("E40").Value = "  -8.48%  "
$ws.Range("D41").Value = "'2.69"
$ws.Range("E41").Value = "  -12.23%  "
$ws.Range("E42").Value = "  -11.18%  "
$ws.Range("D43").Value = "'3.09"
$ws.Range("E43").Value = "  +16.70%  "
$ws.Range("B44").Value = "Monero"
$ws.Range("C44").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D44").Value = "'142.24"
$ws.Range("E44").Value = "  -3.44%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'26.00"
$ws.Range("E45").Value = "  -6.64%  "
$ws.Range("D46").Value = "'1.97"
$ws.Range("E46").Value = "  -5.60%  "
$ws.Range("D47").Value = "'3.09"
$ws.Range("E47").Value = "  -10.05%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'4.09"
$ws.Range("E48").Value = "  -5.57%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.48"
$ws.Range("E49").Value = "  -7.45%  "
$ws.Range("D50").Value = "'2.67"
$ws.Range("E50").Value = "  -8.64%  "
$ws.Range("D51").Value = "'0.279"
$ws.Range("E51").Value = "  -8.69%  "
